$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtering save games) for rows 2-10, columns B-G
$data = @{
    2 = @{ B = 3.230985683306322;    C = 1.667794583268128;    D = 0.8054896365839992;    E = 8.660232485948974;  F = 1; G = 14.36450238910742 }
    3 = @{ B = 3.230985683306322;    C = 1.667794583268128;    D = 0.1575252929769615;    E = 0.496779210170732;  F = 0; G = 5.553084769722144 }
    4 = @{ B = 0.6753301551942219;   C = 1.667794583268128;    D = 3.900430680208489;     E = 0.496779210170732;  F = 1; G = 6.740334628841572 }
    5 = @{ B = 3.230985683306322;    C = 1.667794583268128;    D = 0.8054896365839992;    E = 0.496779210170732;  F = 0; G = 6.201049113329182 }
    6 = @{ B = 1.459612070389937;    C = 1.667794583268128;    D = 0.1575252929769615;    E = 0.496779210170732;  F = 0; G = 3.781711156805759 }
    7 = @{ B = 3.230985683306322;    C = 1.667794583268128;    D = 3.900430680208489;     E = 0.496779210170732;  F = 0; G = 9.295990156953671 }
    8 = @{ B = 0.6753301551942219;   C = 1.667794583268128;    D = 0.8054896365839992;    E = 0.496779210170732;  F = 1; G = 3.645393585217082 }
    9 = @{ B = 0.6753301551942219;   C = 1.667794583268128;    D = 0.8054896365839992;    E = 0.496779210170732;  F = 1; G = 3.645393585217082 }
    10 = @{ B = 0.000002317355952907718; C = 0.00007097389502863649; D = 0.1575252929769615; E = 0.496779210170732;  F = 0; G = 0.6543777943986751 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
